# Daily attendance processing - 2025-10-24 05:44:50
# Rotates the "Recorded By" (column G) comma-separated list so that the
# last entry moves to the front, for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7 ("Recorded By")
    $value = $cell.Value()

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = @($parts[$parts.Count - 1]) + $parts[0..($parts.Count - 2)]
            $cell.Value = [string]::Join(", ", $rotated)
        }
    }
}
